# iter8 illustration.pptx update:
#  - bump cached date fields from 7/21/21 -> 7/28/21 (slide master + every slide layout)
#  - "Safely rollout winning version" -> "Assess and safely rollout winning version of app"
#  - "Use builtin metrics or custom metrics from any DB"
#        -> "Use built-in metrics collector or custom metrics from any DB"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetime" field text wherever it appears (ppPlaceholderDate = 16)
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes, $newDate) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$newDate = "7/28/21"

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes $newDate

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes $newDate
}

# ---------------------------------------------------------------------------
# 2) "Rounded Rectangle 54" callout text
#    "Safely rollout winning version" -> "Assess and safely rollout winning version of app"
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$rr54 = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    if ($slide.Shapes.Item($i).Name -eq "Rounded Rectangle 54") {
        $rr54 = $slide.Shapes.Item($i)
        break
    }
}

$tr = $rr54.TextFrame.TextRange

# Drop the leading "Safely " run entirely (first 7 characters).
$tr.Characters(1, 7).Text = ""

# The remaining text now starts with "rollout winning version"; rewrite the
# "rollout " run (still its own run, 8 chars) in place so it keeps its own
# (dirty="0") run formatting.
$tr.Characters(1, 8).Text = "Assess and safely rollout "

# The bold "winning version" run follows right after; rewrite it in place too.
$prefixLen = "Assess and safely rollout ".Length
$tr.Characters($prefixLen + 1, "winning version".Length).Text = "winning version of app"

# ---------------------------------------------------------------------------
# 3) "Rounded Rectangle 65" callout text
#    "Use builtin metrics or custom metrics from any DB"
#      -> "Use built-in metrics collector or custom metrics from any DB"
# ---------------------------------------------------------------------------
$rr65 = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    if ($slide.Shapes.Item($i).Name -eq "Rounded Rectangle 65") {
        $rr65 = $slide.Shapes.Item($i)
        break
    }
}

$tr2 = $rr65.TextFrame.TextRange

# Runs 2-4 of this paragraph are " ", "builtin", " " (chars 4-12 == " builtin ").
# Collapse them into a single " built-in " run, keeping run 2's formatting.
$builtinReplacement = " built-in "
$tr2.Characters(4, 9).Text = $builtinReplacement

# "metrics or " (now starting right after "Use" + " built-in ") becomes "metrics collector or "
$afterBuiltIn = "Use".Length + $builtinReplacement.Length
$tr2.Characters($afterBuiltIn + 1, "metrics or ".Length).Text = "metrics collector or "
